$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 129.66667
$ws.Range("I33").Value = 129.66667
$ws.Range("K33").Value = 129.66667
$ws.Range("M33").Value = 99.33332999999999
$ws.Range("H113").Value = 5024.5
$ws.Range("I113").Value = 3251
$ws.Range("J113").Value = 5615.6665
$ws.Range("K113").Value = 3251
$ws.Range("L113").Value = 5615.6665
$ws.Range("M113").Value = 3
$ws.Range("N113").Value = -12123.6665
$ws.Range("H116").Value = 12000
$ws.Range("J116").Value = 12667
$ws.Range("L116").Value = 12667
$ws.Range("N116").Value = -19551
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2275349.8
$ws.Range("I32").Value = 2447.5557
$ws.Range("K32").Value = 2447.5557
$ws.Range("M32").Value = -2160.5557
$ws.Range("H45").Value = 3339.9443
$ws.Range("I45").Value = 2850.2727
$ws.Range("K45").Value = 2850.2727
$ws.Range("M45").Value = -2473.2727
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H132").Value = 2692.0625
$ws.Range("I132").Value = 2313.3845
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 6940.1535
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -4410.1535
$ws.Range("N132").Value = -18059
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2168.1428
$ws.Range("I99").Value = 2436
$ws.Range("K99").Value = 2436
$ws.Range("M99").Value = -938
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2105.8333
$ws.Range("I58").Value = 1570.0769
$ws.Range("J58").Value = 3498.8
$ws.Range("K58").Value = 1570.0769
$ws.Range("L58").Value = 3498.8
$ws.Range("M58").Value = -1367.0769
$ws.Range("N58").Value = -3904.8
$ws.Range("H136").Value = 2105.8333
$ws.Range("I136").Value = 1570.0769
$ws.Range("J136").Value = 3498.8
$ws.Range("K136").Value = 4710.2307
$ws.Range("L136").Value = 10496.4
$ws.Range("M136").Value = -2160.2307
$ws.Range("N136").Value = -15596.4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1659.8125
$ws.Range("I14").Value = 1659.8125
$ws.Range("K14").Value = 4979.4375
$ws.Range("M14").Value = -4806.4375
$ws.Range("H34").Value = 933.9375
$ws.Range("J34").Value = 1288.2727
$ws.Range("L34").Value = 3864.8181
$ws.Range("N34").Value = -4032.8181
$ws.Range("H39").Value = 4845.4287
$ws.Range("J39").Value = 5322.316
$ws.Range("L39").Value = 15966.948
$ws.Range("N39").Value = -16554.948
$ws.Range("H40").Value = 297.1111
$ws.Range("J40").Value = 341
$ws.Range("L40").Value = 1364
$ws.Range("N40").Value = -1502
$ws.Range("H55").Value = 11530.75
$ws.Range("J55").Value = 12478
$ws.Range("L55").Value = 37434
$ws.Range("N55").Value = -37788
$ws.Range("H87").Value = 3007
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 3007
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 70.375
$ws.Range("I97").Value = 84.5
$ws.Range("J97").Value = 56.25
$ws.Range("K97").Value = 84.5
$ws.Range("L97").Value = 56.25
$ws.Range("M97").Value = 411.5
$ws.Range("N97").Value = -1048.25
$ws.Range("H113").Value = 4274.615
$ws.Range("I113").Value = 2285.5557
$ws.Range("K113").Value = 2285.5557
$ws.Range("M113").Value = -115.5556999999999
$ws.Range("H132").Value = 2351
$ws.Range("I132").Value = 1510.5454
$ws.Range("K132").Value = 4531.6362
$ws.Range("M132").Value = -2001.6362
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8794.5
$ws.Range("I7").Value = 8792.333000000001
$ws.Range("K7").Value = 8792.333000000001
$ws.Range("M7").Value = -8680.333000000001
$ws.Range("H46").Value = 6417.5
$ws.Range("I46").Value = 920
$ws.Range("J46").Value = 8250
$ws.Range("K46").Value = 920
$ws.Range("L46").Value = 8250
$ws.Range("M46").Value = -732
$ws.Range("N46").Value = -8626
$ws.Range("H61").Value = 5607.1113
$ws.Range("I61").Value = 4344
$ws.Range("J61").Value = 7186
$ws.Range("K61").Value = 4344
$ws.Range("L61").Value = 7186
$ws.Range("M61").Value = -4142
$ws.Range("N61").Value = -7590
$ws.Range("H113").Value = 5607.1113
$ws.Range("I113").Value = 4344
$ws.Range("J113").Value = 7186
$ws.Range("K113").Value = 4344
$ws.Range("L113").Value = 7186
$ws.Range("M113").Value = -2174
$ws.Range("N113").Value = -11526
$ws.Range("H126").Value = 8794.5
$ws.Range("I126").Value = 8792.333000000001
$ws.Range("K126").Value = 26376.999
$ws.Range("M126").Value = -23906.999
$ws.Range("H132").Value = 13221.625
$ws.Range("I132").Value = 7698.25
$ws.Range("K132").Value = 23094.75
$ws.Range("M132").Value = -20564.75
$ws.Range("H136").Value = 3198.818
$ws.Range("I136").Value = 2483.75
$ws.Range("J136").Value = 4056.9
$ws.Range("K136").Value = 7451.25
$ws.Range("L136").Value = 12170.7
$ws.Range("M136").Value = -4901.25
$ws.Range("N136").Value = -17270.7
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 31310.5
$ws.Range("J41").Value = 39955.668
$ws.Range("L41").Value = 39955.668
$ws.Range("N41").Value = -40735.668
$ws.Range("H122").Value = 2666.3103
$ws.Range("I122").Value = 2793.44
$ws.Range("J122").Value = 1871.75
$ws.Range("K122").Value = 8380.32
$ws.Range("L122").Value = 5615.25
$ws.Range("M122").Value = -5930.32
$ws.Range("N122").Value = -10515.25
$ws.Range("H132").Value = 2613
$ws.Range("I132").Value = 2613
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7839
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = -5309
$ws.Range("M132").ClearContents()
